$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.355.77"
$ws.Range("E2").Value = "  -2.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.884.72"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.66"
$ws.Range("E5").Value = "  -4.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.64"
$ws.Range("E6").Value = "  -2.35%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.883.64"
$ws.Range("E9").Value = "  -1.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.90"
$ws.Range("E10").Value = "  -5.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.146"
$ws.Range("E11").Value = "  -1.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.432"
$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.97"
$ws.Range("E14").Value = "  -1.88%  "

$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.354.27"
$ws.Range("E16").Value = "  -2.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.367.03"
$ws.Range("E17").Value = "  -2.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.56"
$ws.Range("E18").Value = "  -1.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.886.63"
$ws.Range("E19").Value = "  -1.89%  "

$ws.Range("E20").Value = "  -1.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.07"
$ws.Range("E21").Value = "  -2.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.655"
$ws.Range("E22").Value = "  -1.36%  "

$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.89"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").Value = "  -9.94%  "

$ws.Range("E28").Value = "  -5.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("E29").Value = "  +3.86%  "

$ws.Range("E30").Value = "  -2.22%  "

$ws.Range("E31").Value = "  -4.04%  "

$ws.Range("E32").Value = "  -7.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.43"
$ws.Range("E35").Value = "  -3.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.957"
$ws.Range("E36").Value = "  -3.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.39"
$ws.Range("E37").Value = "  -3.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.85"
$ws.Range("E38").Value = "  -1.64%  "

$ws.Range("E39").Value = "  -7.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").Value = "  -4.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.20"
$ws.Range("E41").Value = "  -2.98%  "

$ws.Range("E42").Value = "  -2.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.80"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("E44").Value = "  -4.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.710.31"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.05"
$ws.Range("E46").Value = "  -1.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0335"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "342.09"
$ws.Range("E48").Value = "  -5.13%  "

$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("E51").Value = "  -4.75%  "

Write-Host "Applied cryptos update"